# Reproduces the commit "Add files via upload / init commit":
#   - Sheet "OSMEI Model Equations" (2nd tab) is renamed to "Modified Model Equations"
#   - That sheet becomes the active/selected tab (it was the "Variable List" tab before)
#   - The selection/active cell on "Modified Model Equations" moves to G2 (was E167)
#   - The selection/active cell on "Paper Model Equations" moves to F7 (was E182:G182)

$wb = $excel.ActiveWorkbook

$wsVariableList = $wb.Worksheets.Item(1)   # "Variable List"
$wsEquations    = $wb.Worksheets.Item(2)   # "OSMEI Model Equations" -> "Modified Model Equations"
$wsPaperModel   = $wb.Worksheets.Item(3)   # "Paper Model Equations"

# Rename the sheet.
$wsEquations.Name = "Modified Model Equations"

# Update the selection on "Paper Model Equations" first (it stays a background tab).
$wsPaperModel.Activate()
$wsPaperModel.Range("F7").Select() | Out-Null

# Make "Modified Model Equations" the active tab with its new selection - this is the
# last sheet activated, so it ends up as the tabSelected / activeTab sheet on save.
$wsEquations.Activate()
$wsEquations.Range("G2").Select() | Out-Null
